$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.668.14"
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").Value = "3.337.52"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.24"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "183.47"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.98%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.648"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +4.70%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "3.337.51"
$ws.Range("E9").Value = "  +0.43%  "
$ws.Range("E10").Value = "  -0.48%  "
$ws.Range("E11").Value = "  +2.73%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.403"
$ws.Range("D12").ClearFormats()
$ws.Range("D13").Value = "3.921.56"
$ws.Range("E13").Value = "  +0.35%  "
$ws.Range("E14").Value = "  -2.47%  "
$ws.Range("D15").Value = "66.657.79"
$ws.Range("E15").Value = "  +0.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.77"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("E17").Value = "  -0.48%  "
$ws.Range("D18").Value = "3.337.88"
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "427.06"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.63%  "
$ws.Range("E20").Value = "  -1.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.21"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.44"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.48%  "
$ws.Range("E23").Value = "  -1.36%  "
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("E25").Value = "  +0.48%  "
$ws.Range("D26").Value = "3.472.82"
$ws.Range("E26").Value = "  -0.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.515"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.30%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.207"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +7.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000116"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.45%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.09"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.09%  "
$ws.Range("E32").Value = "  -0.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "22.49"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.03%  "
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.24"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.68"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.75%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.20"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "160.46"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.31%  "
$ws.Range("E39").Value = "  -1.02%  "
$ws.Range("E40").Value = "  +1.91%  "
$ws.Range("D41").Value = "2.871.13"
$ws.Range("E41").Value = "  +1.70%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.60"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -4.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.35"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.763"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.81%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "39.80"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.99%  "
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0665"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.41%  "
$ws.Range("E47").Value = "  -2.68%  "
$ws.Range("E48").Value = "  +1.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.35"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "314.67"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0274"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.22%  "
